$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games) for rows 2-11, columns B-G
$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043;  E = 0.5333859586016987; F = 0; G = 5.582307763322248 }
    3  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 18.71679738969934;   E = 0.5333859586016987; F = 1; G = 24.14949828602258 }
    4  = @{ B = 0.6545652718822623; C = 1.626987699542094;   D = 3.223369029078222;   E = 0.5333859586016987; F = 1; G = 6.038307959104277 }
    5  = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.7210945179870265;  E = 0.5333859586016987; F = 0; G = 4.327115817150455 }
    6  = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.7210945179870265;  E = 0.5333859586016987; F = 0; G = 4.327115817150455 }
    7  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265;  E = 0.5333859586016987; F = 1; G = 6.15379541431027 }
    8  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;   E = 0.5333859586016987; F = 0; G = 8.656069925401464 }
    9  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;   E = 0.5333859586016987; F = 0; G = 8.656069925401464 }
    10 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265;  E = 0.5333859586016987; F = 0; G = 6.15379541431027 }
    11 = @{ B = 0.2881169905109251; C = 0.04103571897497393; D = 0.7210945179870265;  E = 0.5333859586016987; F = 1; G = 1.583633186074624 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
}
